$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1900
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1900
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 5700
$ws.Range("N58").Value = -6000
$ws.Range("M58").ClearContents()
$ws.Range("H64").Value = 3999.5
$ws.Range("I64").Value = 3999
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3999
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3751
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 3999.5
$ws.Range("I67").Value = 3999
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3999
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -3141
$ws.Range("N67").Value = -5716
$ws.Range("H112").Value = 4899
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4899
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 14697
$ws.Range("N112").Value = -16913
$ws.Range("H135").Value = 1517.8334
$ws.Range("I135").Value = 1637.75
$ws.Range("J135").Value = 1278
$ws.Range("K135").Value = 14739.75
$ws.Range("L135").Value = 11502
$ws.Range("M135").Value = -12204.75
$ws.Range("N135").Value = -16572
$ws.Range("H137").Value = 2263.625
$ws.Range("I137").Value = 2292.4
$ws.Range("J137").Value = 2215.6667
$ws.Range("K137").Value = 6877.200000000001
$ws.Range("L137").Value = 6647.000100000001
$ws.Range("M137").Value = -4327.200000000001
$ws.Range("N137").Value = -11747.0001
$ws.Range("H138").Value = 3449.75
$ws.Range("I138").Value = 1607.8334
$ws.Range("J138").Value = 3952.0908
$ws.Range("K138").Value = 4823.5002
$ws.Range("L138").Value = 11856.2724
$ws.Range("M138").Value = 316.4997999999996
$ws.Range("N138").Value = -22136.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2959.4736
$ws.Range("I32").Value = 2219.2856
$ws.Range("J32").Value = 11595
$ws.Range("K32").Value = 2219.2856
$ws.Range("L32").Value = 11595
$ws.Range("M32").Value = -1932.2856
$ws.Range("H36").Value = 6513
$ws.Range("I36").Value = 5026
$ws.Range("J36").Value = 8000
$ws.Range("K36").Value = 5026
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = -4680
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4199.2
$ws.Range("I20").Value = 4199.2
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 4199.2
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -3952.2
$ws.Range("H44").Value = 65000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 65000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 65000
$ws.Range("N44").Value = -65994
$ws.Range("H99").Value = 4870.357
$ws.Range("I99").Value = 4976.154
$ws.Range("J99").Value = 3495
$ws.Range("K99").Value = 4976.154
$ws.Range("L99").Value = 3495
$ws.Range("M99").Value = -3478.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5693.4443
$ws.Range("I31").Value = 6373.6665
$ws.Range("J31").Value = 4333
$ws.Range("K31").Value = 6373.6665
$ws.Range("L31").Value = 4333
$ws.Range("M31").Value = -6078.6665
$ws.Range("N31").Value = -4923
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9621
$ws.Range("H34").Value = 5693.4443
$ws.Range("I34").Value = 6373.6665
$ws.Range("J34").Value = 4333
$ws.Range("K34").Value = 6373.6665
$ws.Range("L34").Value = 4333
$ws.Range("M34").Value = -6171.6665
$ws.Range("N34").Value = -4737
$ws.Range("H60").Value = 28125
$ws.Range("I60").Value = 22000
$ws.Range("J60").Value = 29000
$ws.Range("K60").Value = 22000
$ws.Range("L60").Value = 29000
$ws.Range("M60").Value = -21489
$ws.Range("N60").Value = -30022
$ws.Range("H62").Value = 8046.75
$ws.Range("I62").Value = 9339
$ws.Range("J62").Value = 5893
$ws.Range("K62").Value = 9339
$ws.Range("L62").Value = 5893
$ws.Range("M62").Value = -8715
$ws.Range("N62").Value = -7141
$ws.Range("H65").Value = 8046.75
$ws.Range("I65").Value = 9339
$ws.Range("J65").Value = 5893
$ws.Range("K65").Value = 46695
$ws.Range("L65").Value = 29465
$ws.Range("M65").Value = -43575
$ws.Range("N65").Value = -35705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 31.5
$ws.Range("I12").Value = 100.5
$ws.Range("J12").Value = 21.642857
$ws.Range("K12").Value = 301.5
$ws.Range("L12").Value = 64.92857100000001
$ws.Range("M12").Value = -128.5
$ws.Range("N12").Value = -410.928571
$ws.Range("H49").Value = 999
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 999
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 2997
$ws.Range("N49").Value = -3309
$ws.Range("H109").Value = 2699
$ws.Range("I109").Value = 2699
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 8097
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -7057
$ws.Range("H111").Value = 2999
$ws.Range("I111").Value = 2999
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 8997
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -5930
$ws.Range("H114").Value = 2214
$ws.Range("I114").Value = 766.3333
$ws.Range("J114").Value = 3299.75
$ws.Range("K114").Value = 2298.9999
$ws.Range("L114").Value = 9899.25
$ws.Range("M114").Value = 955.0001000000002
$ws.Range("N114").Value = -16407.25
$ws.Range("H122").Value = 2876.25
$ws.Range("I122").Value = 2702
$ws.Range("J122").Value = 2934.3333
$ws.Range("K122").Value = 24318
$ws.Range("L122").Value = 26408.9997
$ws.Range("M122").Value = -21868
$ws.Range("N122").Value = -31308.9997
$ws.Range("H132").Value = 1971.7142
$ws.Range("I132").Value = 1575.75
$ws.Range("J132").Value = 2499.6667
$ws.Range("K132").Value = 14181.75
$ws.Range("L132").Value = 22497.0003
$ws.Range("M132").Value = -11651.75
$ws.Range("N132").Value = -27557.0003
$ws.Range("H141").Value = 1899.6
$ws.Range("I141").Value = 1899.6
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5698.799999999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -518.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5001.5
$ws.Range("I70").Value = 5003.5
$ws.Range("J70").Value = 4999.5
$ws.Range("K70").Value = 5003.5
$ws.Range("L70").Value = 4999.5
$ws.Range("M70").Value = -4733.5
$ws.Range("N70").Value = -5539.5
$ws.Range("H73").Value = 5001.5
$ws.Range("I73").Value = 5003.5
$ws.Range("J73").Value = 4999.5
$ws.Range("K73").Value = 5003.5
$ws.Range("L73").Value = 4999.5
$ws.Range("M73").Value = -4067.5
$ws.Range("N73").Value = -6871.5
$ws.Range("H122").Value = 2812.3333
$ws.Range("I122").Value = 2972
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 8916
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -6466
$ws.Range("N122").Value = -11899.9999
$ws.Range("H126").Value = 6749.5
$ws.Range("I126").Value = 6749.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 20248.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17778.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 6089.7144
$ws.Range("I132").Value = 6210.6665
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 18631.9995
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -16101.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 17256
$ws.Range("I43").Value = 17256
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 17256
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -17063
$ws.Range("H122").Value = 3575.111
$ws.Range("I122").Value = 3575.111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10725.333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8275.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 2000
$ws.Range("N20").Value = -2480
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H34").Value = 18000
$ws.Range("I34").Value = 18000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -17797
$ws.Range("H37").Value = 15410
$ws.Range("I37").Value = 9017.333000000001
$ws.Range("J37").Value = 24999
$ws.Range("K37").Value = 9017.333000000001
$ws.Range("L37").Value = 24999
$ws.Range("M37").Value = -8814.333000000001
$ws.Range("N37").Value = -25405
$ws.Range("H46").Value = 9429
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 9429
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 9429
$ws.Range("N46").Value = -9891
$ws.Range("H122").Value = 2050.7727
$ws.Range("I122").Value = 1690.421
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 5071.263
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -2621.263
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3160.6316
$ws.Range("I132").Value = 1440.9375
$ws.Range("J132").Value = 12332.333
$ws.Range("K132").Value = 4322.8125
$ws.Range("L132").Value = 36996.999
$ws.Range("M132").Value = -1792.8125
$ws.Range("N132").Value = -42056.999
$ws.Range("H134").Value = 9429
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 9429
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 28287
$ws.Range("N134").Value = -33357
$ws.Range("H136").Value = 3068.2307
$ws.Range("I136").Value = 2535.4546
$ws.Range("J136").Value = 5998.5
$ws.Range("K136").Value = 7606.3638
$ws.Range("L136").Value = 17995.5
$ws.Range("M136").Value = -5056.3638
$ws.Range("N136").Value = -23095.5

